$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Create the new "2022-Q4" sheet by duplicating "2022-Q3" (so it keeps
#    the exact same header/column styling) right before it, then trim it
#    down to the header + 2 fund rows and overwrite the values.
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item(2)
$q3.Copy($q3)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# Drop the extra data rows (the Q3 sheet has 7 fund rows, Q4 only has 2).
$q4.Range("A4:H8").EntireRow.Delete()

# Columns B-G are plain text in the source data, so force a text
# number-format before writing so Excel doesn't auto-convert numeric-
# looking strings ("161039", "20.17", ...) into real numbers.
$q4.Range("B2:G3").NumberFormat = "@"

$q4.Range("A2").Value = 0
$q4.Range("B2").Value = "161039"
$q4.Range("C2").Value = "富国中证1000指数增强（LOF）A"
$q4.Range("D2").Value = "20.17"
$q4.Range("E2").Value = "91.85"
$q4.Range("F2").Value = "0.62"
$q4.Range("G2").Value = "0.1251"
$q4.Range("H2").Value = 9

$q4.Range("A3").Value = 1
$q4.Range("B3").Value = "013331"
$q4.Range("C3").Value = "富国中证1000指数增强（LOF）C"
$q4.Range("D3").Value = "6.90"
$q4.Range("E3").Value = "91.85"
$q4.Range("F3").Value = "0.62"
$q4.Range("G3").Value = "0.0428"
$q4.Range("H3").Value = 9

# The text number-format was only needed to land the values as text; put
# the cell style back the way the untouched source cells had it (no
# explicit style override).
$q4.Range("B2:G3").Style = "Normal"

# ---------------------------------------------------------------------
# 2) Add the matching summary row to the "总计" sheet: insert a new row 2
#    (pushing the existing quarters down one) and fill it in with the
#    2022-Q4 totals.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)
$total.Rows.Item(2).Insert()

$total.Range("A2").Value = 0
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Range("B2:D2").Style = "Normal"
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.17

# Renumber the index column (A) for the rows that shifted down so it
# keeps counting 0,1,2,3,4,5 from top to bottom.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5
